# Updated symbol list on Sat Dec 31 13:56:18 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as literal text
# (inline/shared strings, not numbers) so that exact textual formatting -
# including trailing zeros and "long-tail" decimals - survives untouched.
# Excel's normal Range.Value assignment auto-detects numeric-looking
# strings and coerces them into real floating point numbers, which both
# loses trailing zeros (e.g. "0.8130" -> 0.813) and can introduce binary
# floating point artifacts on save. To avoid that, the cell is temporarily
# switched to the Text ("@") number format before the value is written,
# then the formatting is cleared again so the cell ends up back on the
# sheet's default (unstyled) look, matching the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    $rng = $ws.Range($a1)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

function Set-PlainValue($a1, $value) {
    $ws.Range($a1).Value = $value
}

# --- Row 2 (BNB) ---
Set-TextValue "D2" "246.75"

# --- Row 3 (OKB) ---
Set-TextValue "D3" "26.35"

# --- Row 4 (HuobiToken) ---
Set-TextValue "D4" "5.075"

# --- Row 5 (Cronos) ---
Set-TextValue "D5" "0.05606"

# --- Row 6 (KuCoinToken) ---
Set-TextValue "D6" "6.505"

# --- Row 8 (MXToken) ---
Set-TextValue "D8" "0.8130"

# --- Row 9 (FTXToken) ---
Set-TextValue "D9" "0.8412"

# --- Row 10: One -> WazirX ---
Set-PlainValue "B10" "WazirX"
Set-PlainValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue  "D10" "0.1345"
Set-PlainValue "E10" "9WazirXWRX"

# --- Row 11: WazirX -> BitrueCoin ---
Set-PlainValue "B11" "BitrueCoin"
Set-PlainValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue  "D11" "0.02826"
Set-PlainValue "E11" "10BitrueCoinBTR"

# --- Row 12: BitrueCoin -> BitMartToken ---
Set-PlainValue "B12" "BitMartToken"
Set-PlainValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue  "D12" "0.09397"
Set-PlainValue "E12" "11BitMartTokenBMX"

# --- Row 13: BitMartToken -> BitForexToken ---
Set-PlainValue "B13" "BitForexToken"
Set-PlainValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue  "D13" "0.001527"
Set-PlainValue "E13" "12BitForexTokenBF"

# --- Row 14: BitForexToken -> One ---
Set-PlainValue "B14" "One"
Set-PlainValue "C14" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue  "D14" "0.0005987"
Set-PlainValue "E14" "13OneONE"

# --- Row 15 (TigerCash) ---
Set-TextValue "D15" "0.006128"

# --- Row 16 (LEO) ---
Set-TextValue "D16" "3.553"

# --- Row 19 (MandalaExchangeToken) ---
Set-TextValue "D19" "0.06986"

# --- Row 20 (LiechtensteinCryptoassetsExchange) ---
Set-TextValue "D20" "0.03146"

# --- Row 22 (MCDex) ---
Set-TextValue "D22" "3.738"

# --- Row 25 (BitKan) ---
Set-TextValue "D25" "0.001247"

# --- Row 26 (HotbitToken) ---
Set-TextValue "D26" "0.004615"

# --- Row 27 (NitroEx) ---
Set-TextValue "D27" "0.00009597"

# --- Row 28 (UpBots) -- only the "Best in 24h" tag moves here ---
Set-PlainValue "E28" "27UpBotsUBXTBestin24h"

# --- Row 40 (IDEX) ---
Set-TextValue "D40" "0.03665"

# --- Row 41: KickToken -> BKEXToken ---
Set-PlainValue "B41" "BKEXToken"
Set-PlainValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue  "D41" "0.1349"
Set-PlainValue "E41" "40BKEXTokenBKK"

# --- Row 42: BKEXToken -> CEJI ---
Set-PlainValue "B42" "CEJI"
Set-PlainValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue  "D42" "0.002659"
Set-PlainValue "E42" "41CEJICEJI"

# --- Row 43: CEJI -> KickToken ---
Set-PlainValue "B43" "KickToken"
Set-PlainValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue  "D43" "0.003416"
Set-PlainValue "E43" "42KickTokenKICKWorstin24h"

# --- Row 44 (LocalTraders) ---
Set-TextValue "D44" "0.008879"

# --- Row 45 (CoinLion) ---
Set-TextValue "D45" "0.00005287"
